$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update F27:F38 with new classification labels (age cohorts / "ca. YYYY") ---
# The order below matches the order these new values were first entered,
# which determines the order they land in the shared-strings table.
$ws.Range("F27").Value = "1999-2003"

$ws.Range("F38").Value = "ca. 2010"
$ws.Range("F38").ClearFormats()

$ws.Range("F33").Value = "ca. 2005"
$ws.Range("F33").ClearFormats()

$ws.Range("F28").Value = "ca. 2000"
$ws.Range("F28").ClearFormats()

$ws.Range("F29").Value = "ca. 2001"
$ws.Range("F29").ClearFormats()

$ws.Range("F30").Value = "ca. 2002"
$ws.Range("F30").ClearFormats()

$ws.Range("F31").Value = "ca. 2003"
$ws.Range("F31").ClearFormats()

$ws.Range("F32").Value = "ca. 2004"
$ws.Range("F32").ClearFormats()

$ws.Range("F34").Value = "ca. 2006"
$ws.Range("F34").ClearFormats()

$ws.Range("F35").Value = "ca. 2007"
$ws.Range("F35").ClearFormats()

$ws.Range("F36").Value = "ca. 2008"
$ws.Range("F36").ClearFormats()

$ws.Range("F37").Value = "ca. 2009"
$ws.Range("F37").ClearFormats()

# --- Append 100 more rows (102-201) for uploading 10 more datasets, replicating
#     the row 101 pattern: A = row-1, B = 14, C = 0, F = "reserved" ---
# Insert the rows first so the inherited formatting (style s="2"/"3") carries
# down from row 101, matching what Excel does when extending a formatted list.
$ws.Rows("102:201").Insert()

for ($r = 102; $r -le 201; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
    $ws.Cells.Item($r, 2).Value = 14
    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 6).Value = "reserved"
}

# --- Update the view state to match: scrolled so row 13 is at top, F40 selected ---
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("F40").Select()
